{"js": "/* Replace the 100 math-problem cell texts in the table, in row-major\n * order (row 1..20, column 1..5), per the old->new mapping captured\n * from the authoritative diff. Cells are addressed via the table's\n * rows/cells collections so the edit doesn't depend on how the body's\n * flattened paragraph list happens to be laid out.\n */\nconst mapping = [[\"8+81=\", \"19-17=\"], [\"14+33=\", \"89-84=\"], [\"50+20=\", \"46+49=\"], [\"30+48=\", \"56-23=\"], [\"1+90=\", \"44+17=\"], [\"90-75=\", \"77-13=\"], [\"95-68=\", \"65-45=\"], [\"32+32=\", \"81-50=\"], [\"64-3=\", \"21+28=\"], [\"23+68=\", \"65+21=\"], [\"90-24=\", \"73+22=\"], [\"0+46=\", \"98-37=\"], [\"15+34=\", \"88-87=\"], [\"51-27=\", \"89-72=\"], [\"34+8=\", \"64-18=\"], [\"27-15=\", \"83-51=\"], [\"42+32=\", \"90-28=\"], [\"82-55=\", \"20-6=\"], [\"83-31=\", \"35+39=\"], [\"52-44=\", \"33+22=\"], [\"49+16=\", \"11+86=\"], [\"11+15=\", \"59+37=\"], [\"33-22=\", \"73+16=\"], [\"50-22=\", \"56+39=\"], [\"40-29=\", \"4+21=\"], [\"5+52=\", \"84-3=\"], [\"2+46=\", \"93-91=\"], [\"57-9=\", \"59-25=\"], [\"83-62=\", \"78-39=\"], [\"5-4=\", \"3+67=\"], [\"32-1=\", \"52+35=\"], [\"20+77=\", \"38+60=\"], [\"30-4=\", \"6+79=\"], [\"74-21=\", \"74+24=\"], [\"42+53=\", \"44+25=\"], [\"81-73=\", \"42-31=\"], [\"58+16=\", \"44+39=\"], [\"19+27=\", \"2+36=\"], [\"22+2=\", \"75+19=\"], [\"94-6=\", \"82-29=\"], [\"39+27=\", \"18-15=\"], [\"74-73=\", \"41+31=\"], [\"13+28=\", \"58+15=\"], [\"89-85=\", \"96-77=\"], [\"17+51=\", \"81-38=\"], [\"55-3=\", \"92-24=\"], [\"26+72=\", \"28+47=\"], [\"19+67=\", \"17+0=\"], [\"13+53=\", \"51-2=\"], [\"53-7=\", \"62-55=\"], [\"85-83=\", \"40-6=\"], [\"86-26=\", \"93-36=\"], [\"92-48=\", \"48-21=\"], [\"56+33=\", \"33+52=\"], [\"76-58=\", \"90-19=\"], [\"95-80=\", \"34+8=\"], [\"30+66=\", \"63-38=\"], [\"41-32=\", \"1+69=\"], [\"43-42=\", \"32+42=\"], [\"86-78=\", \"14+21=\"], [\"63+30=\", \"4+57=\"], [\"68+2=\", \"4+21=\"], [\"81-54=\", \"55-5=\"], [\"1+52=\", \"34+5=\"], [\"33+8=\", \"59-7=\"], [\"10+59=\", \"71-8=\"], [\"2+21=\", \"44+30=\"], [\"92-28=\", \"90-8=\"], [\"30+7=\", \"57-48=\"], [\"67-32=\", \"78-66=\"], [\"91-42=\", \"39-24=\"], [\"89-9=\", \"76+7=\"], [\"92-65=\", \"25+19=\"], [\"22+47=\", \"16+23=\"], [\"5+55=\", \"62-16=\"], [\"19+50=\", \"32-25=\"], [\"70-60=\", \"3+75=\"], [\"83-63=\", \"23-21=\"], [\"94-80=\", \"54+10=\"], [\"60+4=\", \"99-45=\"], [\"42-25=\", \"53-38=\"], [\"39+59=\", \"85+8=\"], [\"69-1=\", \"35+8=\"], [\"79-62=\", \"96-58=\"], [\"15+55=\", \"77-28=\"], [\"16-14=\", \"27+20=\"], [\"52+16=\", \"37-19=\"], [\"53+0=\", \"11+9=\"], [\"15+81=\", \"50-16=\"], [\"79-32=\", \"87-30=\"], [\"57+36=\", \"75-49=\"], [\"43-40=\", \"56-32=\"], [\"72-33=\", \"46+15=\"], [\"10+77=\", \"5+4=\"], [\"64+22=\", \"36+0=\"], [\"85-14=\", \"2+96=\"], [\"85+7=\", \"16-12=\"], [\"48+26=\", \"92-10=\"], [\"96-33=\", \"34-14=\"], [\"83-21=\", \"97-50=\"]];\n\nconst COLS = 5;\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length < 1) {\n  throw new Error(\"Expected at least one table in the document body.\");\n}\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nif (rows.items.length * COLS !== mapping.length) {\n  throw new Error(\n    \"Unexpected table size: expected \" + mapping.length +\n    \" cells (\" + COLS + \" per row), found \" + rows.items.length +\n    \" rows\"\n  );\n}\n\n// Load every row's cells up front, then every cell's value, so we only\n// need two sync() round-trips for the whole table.\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    cell.load(\"value\");\n  }\n}\nawait context.sync();\n\nlet i = 0;\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    const [oldText, newText] = mapping[i];\n    const current = cell.value;\n    if (current !== oldText) {\n      throw new Error(\n        \"Mismatch at cell index \" + i + \": expected '\" + oldText +\n        \"' but found '\" + current + \"'\"\n      );\n    }\n    cell.value = newText;\n    i++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 100 math-problem cell texts in the table, in row-major\n# order (row 1..20, column 1..5), per the old->new mapping captured\n# from the authoritative diff. Values are addressed via Table.Cell(r,c)\n# so we don't have to reason about Word's paragraph/cell-mark counting.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$mapping = @(\n  @(\"8+81=\", \"19-17=\"),\n  @(\"14+33=\", \"89-84=\"),\n  @(\"50+20=\", \"46+49=\"),\n  @(\"30+48=\", \"56-23=\"),\n  @(\"1+90=\", \"44+17=\"),\n  @(\"90-75=\", \"77-13=\"),\n  @(\"95-68=\", \"65-45=\"),\n  @(\"32+32=\", \"81-50=\"),\n  @(\"64-3=\", \"21+28=\"),\n  @(\"23+68=\", \"65+21=\"),\n  @(\"90-24=\", \"73+22=\"),\n  @(\"0+46=\", \"98-37=\"),\n  @(\"15+34=\", \"88-87=\"),\n  @(\"51-27=\", \"89-72=\"),\n  @(\"34+8=\", \"64-18=\"),\n  @(\"27-15=\", \"83-51=\"),\n  @(\"42+32=\", \"90-28=\"),\n  @(\"82-55=\", \"20-6=\"),\n  @(\"83-31=\", \"35+39=\"),\n  @(\"52-44=\", \"33+22=\"),\n  @(\"49+16=\", \"11+86=\"),\n  @(\"11+15=\", \"59+37=\"),\n  @(\"33-22=\", \"73+16=\"),\n  @(\"50-22=\", \"56+39=\"),\n  @(\"40-29=\", \"4+21=\"),\n  @(\"5+52=\", \"84-3=\"),\n  @(\"2+46=\", \"93-91=\"),\n  @(\"57-9=\", \"59-25=\"),\n  @(\"83-62=\", \"78-39=\"),\n  @(\"5-4=\", \"3+67=\"),\n  @(\"32-1=\", \"52+35=\"),\n  @(\"20+77=\", \"38+60=\"),\n  @(\"30-4=\", \"6+79=\"),\n  @(\"74-21=\", \"74+24=\"),\n  @(\"42+53=\", \"44+25=\"),\n  @(\"81-73=\", \"42-31=\"),\n  @(\"58+16=\", \"44+39=\"),\n  @(\"19+27=\", \"2+36=\"),\n  @(\"22+2=\", \"75+19=\"),\n  @(\"94-6=\", \"82-29=\"),\n  @(\"39+27=\", \"18-15=\"),\n  @(\"74-73=\", \"41+31=\"),\n  @(\"13+28=\", \"58+15=\"),\n  @(\"89-85=\", \"96-77=\"),\n  @(\"17+51=\", \"81-38=\"),\n  @(\"55-3=\", \"92-24=\"),\n  @(\"26+72=\", \"28+47=\"),\n  @(\"19+67=\", \"17+0=\"),\n  @(\"13+53=\", \"51-2=\"),\n  @(\"53-7=\", \"62-55=\"),\n  @(\"85-83=\", \"40-6=\"),\n  @(\"86-26=\", \"93-36=\"),\n  @(\"92-48=\", \"48-21=\"),\n  @(\"56+33=\", \"33+52=\"),\n  @(\"76-58=\", \"90-19=\"),\n  @(\"95-80=\", \"34+8=\"),\n  @(\"30+66=\", \"63-38=\"),\n  @(\"41-32=\", \"1+69=\"),\n  @(\"43-42=\", \"32+42=\"),\n  @(\"86-78=\", \"14+21=\"),\n  @(\"63+30=\", \"4+57=\"),\n  @(\"68+2=\", \"4+21=\"),\n  @(\"81-54=\", \"55-5=\"),\n  @(\"1+52=\", \"34+5=\"),\n  @(\"33+8=\", \"59-7=\"),\n  @(\"10+59=\", \"71-8=\"),\n  @(\"2+21=\", \"44+30=\"),\n  @(\"92-28=\", \"90-8=\"),\n  @(\"30+7=\", \"57-48=\"),\n  @(\"67-32=\", \"78-66=\"),\n  @(\"91-42=\", \"39-24=\"),\n  @(\"89-9=\", \"76+7=\"),\n  @(\"92-65=\", \"25+19=\"),\n  @(\"22+47=\", \"16+23=\"),\n  @(\"5+55=\", \"62-16=\"),\n  @(\"19+50=\", \"32-25=\"),\n  @(\"70-60=\", \"3+75=\"),\n  @(\"83-63=\", \"23-21=\"),\n  @(\"94-80=\", \"54+10=\"),\n  @(\"60+4=\", \"99-45=\"),\n  @(\"42-25=\", \"53-38=\"),\n  @(\"39+59=\", \"85+8=\"),\n  @(\"69-1=\", \"35+8=\"),\n  @(\"79-62=\", \"96-58=\"),\n  @(\"15+55=\", \"77-28=\"),\n  @(\"16-14=\", \"27+20=\"),\n  @(\"52+16=\", \"37-19=\"),\n  @(\"53+0=\", \"11+9=\"),\n  @(\"15+81=\", \"50-16=\"),\n  @(\"79-32=\", \"87-30=\"),\n  @(\"57+36=\", \"75-49=\"),\n  @(\"43-40=\", \"56-32=\"),\n  @(\"72-33=\", \"46+15=\"),\n  @(\"10+77=\", \"5+4=\"),\n  @(\"64+22=\", \"36+0=\"),\n  @(\"85-14=\", \"2+96=\"),\n  @(\"85+7=\", \"16-12=\"),\n  @(\"48+26=\", \"92-10=\"),\n  @(\"96-33=\", \"34-14=\"),\n  @(\"83-21=\", \"97-50=\")\n)\n\n$cols = 5\nfor ($i = 0; $i -lt $mapping.Count; $i++) {\n  $row = [int]([math]::Floor($i / $cols)) + 1\n  $col = ($i % $cols) + 1\n  $oldText = $mapping[$i][0]\n  $newText = $mapping[$i][1]\n  $cell = $t.Cell($row, $col)\n  $rng = $cell.Range\n  # Cell.Range.Text carries the trailing end-of-cell mark (\"`r`a\"); strip it\n  # before comparing against the plain cell text captured in $mapping.\n  $current = $rng.Text.TrimEnd([char]13, [char]7)\n  if ($current -ne $oldText) {\n    throw \"Mismatch at row=$row col=${col}: expected [$oldText] but found [$current]\"\n  }\n  $rng.Text = $newText\n}\n\n"}
